# Slide 39: the table's merged header cell (row 1, col 1) currently reads
# "Main Sources of Spatially Referenced Data". Trim the leading "Main "
# so it reads "Sources of Spatially Referenced Data", preserving all
# existing run formatting (bold, size, color, etc.).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(39)
$shp = $s.Shapes.Item(1)
$tbl = $shp.Table
$cell = $tbl.Cell(1, 1)
$cell.Shape.TextFrame.TextRange.Text = "Sources of Spatially Referenced Data"
